# Sprint 2 burndown sheet: update "today"'s remaining-hours entries
# (column C, the "Hours completed" / burn-down input column) for rows
# 4,5,7,8,11,13 and fill in the previously-empty row 14, matching the
# day-by-day numbers that were actually logged. Column E/F are formulas
# and recalc automatically; the scatter-chart's cached values track the
# new F column once Excel recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value  = 0
$ws.Range("C5").Value  = 2
$ws.Range("C7").Value  = 0
$ws.Range("C8").Value  = 2
$ws.Range("C11").Value = 3
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1

# Force a full recalculation so E/F formula caches (and the chart's
# cached series data that reads off them) are refreshed before save.
$excel.CalculateFullRebuild()

# Best-effort: nudge the embedded chart to re-pull its cached point data
# from the now-updated F column (no-ops harmlessly if unsupported).
try {
    $chart = $ws.ChartObjects(1).Chart
    $chart.SetSourceData($ws.Range("A1:A16,E1:F16"))
    $null = $chart.Refresh()
} catch {
}

# Matches the author's last edit being the newly-entered C14.
$null = $ws.Range("C14").Select()
